$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated/consolidated holiday dates for rows 9-57 (A9:A57).
# Rows 2-8 are unchanged. The list below replaces the old values,
# effectively dropping a handful of duplicate Dec 25/Dec 26 entries
# and shifting everything else up.
$newDates = @(
    41639, 41747, 41749, 41750, 41788, 41798, 41799, 41915, 42004,
    42097, 42099, 42100, 42138, 42148, 42149, 42280, 42369, 42454,
    42456, 42457, 42495, 42505, 42506, 42646, 42735, 42839, 42841,
    42842, 42880, 42890, 42891, 43011, 43100, 43189, 43191, 43192,
    43230, 43240, 43241, 43376, 43465, 43574, 43576, 43577, 43615,
    43625, 43626, 43741, 43830
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = 9 + $i
    $ws.Cells.Item($row, 1).Value2 = $newDates[$i]
}

# The rows that used to hold the trailing old dates/flags (58-71) are now
# blank (only the date-formatted style remains on column A).
$ws.Range("A58:B71").ClearContents()

# The sheet used to keep empty placeholder rows out to row 307; now it only
# goes to row 293, so delete the extra trailing rows entirely.
$ws.Range("A294:A307").EntireRow.Delete()

# Update the view state to match where the author left the cursor/scroll.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("A57:XFD58").Select()

Write-Output "edit complete"
